# 无人机车辆协同配送数据.xlsx — add a "车数量" (vehicle count) parameter row
# on the "参数" (parameters) sheet, just below "无人机数量", and nudge the
# saved cursor positions on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 点数据
$ws2 = $wb.Worksheets.Item(2)   # 参数

# Insert a new row above the old row 3 ("无人机单位距离运输成本"), shifting
# everything from row 3 down by one (rows 3-20 -> 4-21).
$ws2.Rows.Item(3).Insert()

# Fill in the new row with the "车数量" parameter and its value.
$ws2.Range("A3").Value = "车数量"
$ws2.Range("B3").Value = 1
$ws2.Rows.Item(3).RowHeight = 25.5

# Restore the sheet selections to match where the author last clicked.
$ws2.Activate() | Out-Null
$ws2.Range("C3").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F6").Select() | Out-Null
